$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.47
$ws.Range("D3").Value = 0.17
$ws.Range("D4").Value = 0.09
$ws.Range("D6").Value = 0.24
$ws.Range("D7").Value = 0.11
$ws.Range("D8").Value = 0.27
$ws.Range("D9").Value = 0.41
$ws.Range("D11").Value = 0.5600000000000001
$ws.Range("D12").Value = 1.48
$ws.Range("D13").Value = 0.11
$ws.Range("D15").Value = 0.16
$ws.Range("D16").Value = 0.27
$ws.Range("D17").Value = 0.6
$ws.Range("D18").Value = 5.13
$ws.Range("D19").Value = 0.54
$ws.Range("D20").Value = 0.33
$ws.Range("D21").Value = 0.61
$ws.Range("D22").Value = 0.1
